$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The original table (rows 8..15) holds the "extr1".."extr8" records.
# Two new records ("line7", "line8") are inserted right after "line6"
# (which sits in row 7), so everything that used to live in rows 8..15
# needs to move down to rows 10..17 first. Shifting via Copy (instead of
# EntireRow.Insert) keeps the existing style table untouched - Insert()
# tends to mint a transient extra cellXf that survives even after the
# cells are overwritten.
# ---------------------------------------------------------------------------
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    $ws.Range("A$r`:E$r").Copy($ws.Range("A$dest`:E$dest"))
}

# Re-sequence column A (0-based running index) for every data row now that
# the table has grown from 14 to 16 data rows (rows 2..17).
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# New row 8: line7
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Updated values for the (now shifted) extr1..extr8 rows
# extr1 -> row 10
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# extr2 -> row 11
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# extr3 -> row 12
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

# extr4 -> row 13
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# extr5 -> row 14
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# extr6 -> row 15
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# extr7 -> row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# extr8 -> row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
